$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.655.68'
$ws.Range("E2").Value = '  +1.26%  '

# Row 3
$ws.Range("D3").Value = '1.559.98'
$ws.Range("E3").Value = '  -0.81%  '

# Row 4
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '210.23'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '0.486'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.42%  '

# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '24.75'
$r.Style = "Normal"
$ws.Range("E8").Value = '  +4.31%  '

# Row 9
$ws.Range("E9").Value = '  -0.14%  '

# Row 10
$ws.Range("E10").Value = '  -0.52%  '

# Row 11
$ws.Range("E11").Value = '  -0.19%  '

# Row 12
$ws.Range("D12").Value = '1.783.11'
$ws.Range("E12").Value = '  -0.78%  '

# Row 13
$ws.Range("D13").Value = '1.563.84'
$ws.Range("E13").Value = '  -0.59%  '

# Row 14
$ws.Range("D14").Value = '28.671.57'
$ws.Range("E14").Value = '  +1.24%  '

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.515'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -0.16%  '

# Row 16
$ws.Range("E16").Value = '  -1.38%  '

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '61.48'
$r.Style = "Normal"
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '230.05'
$r.Style = "Normal"
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '7.37'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -0.61%  '

# Row 20
$ws.Range("E20").Value = '  -1.76%  '

# Row 21
$ws.Range("E21").Value = '  -0.25%  '

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '3.91'
$r.Style = "Normal"
$ws.Range("E22").Value = '  -1.00%  '

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '9.00'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -0.42%  '

# Row 24
$ws.Range("E24").Value = '  +1.57%  '

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '151.09'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '

# Row 26
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '14.76'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '

# Row 27
$ws.Range("E27").Value = '  -0.12%  '

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '

# Row 29
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '6.22'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -2.08%  '

# Row 30
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '0.0460'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -3.97%  '

# Row 31
$ws.Range("E31").Value = '  -2.19%  '

# Row 32
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '3.16'
$r.Style = "Normal"
$ws.Range("E32").Value = '  -0.96%  '

# Row 33
$ws.Range("D33").Value = '1.392.03'
$ws.Range("E33").Value = '  +0.82%  '

# Row 34
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '2.98'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -3.09%  '

# Row 35
$ws.Range("E35").Value = '  -2.75%  '

# Row 36
$ws.Range("E36").Value = '  -1.74%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.29'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -3.11%  '

# Row 38
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '2.65'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +0.52%  '

# Row 39
$ws.Range("E39").Value = '  -0.77%  '

# Row 40
$ws.Range("E40").Value = '  +3.56%  '

# Row 41
$ws.Range("E41").Value = '  -0.33%  '

# Row 42
$ws.Range("E42").Value = '  -0.18%  '

# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.776'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -1.04%  '

# Row 44
$ws.Range("E44").Value = '  +0.39%  '

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '63.93'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +2.69%  '

# Row 46
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '5.27'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -2.09%  '

# Row 47
$ws.Range("D47").Value = '1.695.84'
$ws.Range("E47").Value = '  -0.82%  '

# Row 49
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '85.24'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '

# Row 50
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '43.43'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +5.00%  '

# Row 51
$ws.Range("E51").Value = '  +0.01%  '
